$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New header cell for the next day, matching the formatting of the
# preceding header cell (K1).
$ws.Range("L1").Value = "25-jun"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)

# New daily price values for column L (rows 2-25).
$ws.Range("L2").Value = 84.5
$ws.Range("L3").Value = 74.31
$ws.Range("L4").Value = 75.51000000000001
$ws.Range("L5").Value = 74.94
$ws.Range("L6").Value = 76.39
$ws.Range("L7").Value = 83.45
$ws.Range("L8").Value = 103.14
$ws.Range("L9").Value = 109.54
$ws.Range("L10").Value = 104.74
$ws.Range("L11").Value = 85.23
$ws.Range("L12").Value = 42.63
$ws.Range("L13").Value = 27.98
$ws.Range("L14").Value = 15
$ws.Range("L15").Value = 5.79
$ws.Range("L16").Value = 7.94
$ws.Range("L17").Value = 13.4
$ws.Range("L18").Value = 27.13
$ws.Range("L19").Value = 87.09
$ws.Range("L20").Value = 107.05
$ws.Range("L21").Value = 135
$ws.Range("L22").Value = 143.33
$ws.Range("L23").Value = 141.89
$ws.Range("L24").Value = 135
$ws.Range("L25").Value = 102.26
